$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1081
$ws.Range("I106").Value = 658
$ws.Range("J106").Value = 2350
$ws.Range("K106").Value = 658
$ws.Range("L106").Value = 2350
$ws.Range("M106").Value = -27
$ws.Range("N106").Value = -3612
$ws.Range("H129").Value = 3192.4614
$ws.Range("I129").Value = 2289.7778
$ws.Range("K129").Value = 6869.3334
$ws.Range("M129").Value = -1869.3334
$ws.Range("H132").Value = 68468.53999999999
$ws.Range("I132").Value = 82207.69
$ws.Range("J132").Value = 2062.6667
$ws.Range("K132").Value = 246623.07
$ws.Range("L132").Value = 6188.000100000001
$ws.Range("M132").Value = -244093.07
$ws.Range("N132").Value = -11248.0001
$ws.Range("H138").Value = 1567.0333
$ws.Range("I138").Value = 829.2406999999999
$ws.Range("J138").Value = 2673.7222
$ws.Range("K138").Value = 2487.7221
$ws.Range("L138").Value = 8021.1666
$ws.Range("M138").Value = 2652.2779
$ws.Range("N138").Value = -18301.1666
$ws.Range("H141").Value = 751.3333
$ws.Range("I141").Value = 609
$ws.Range("K141").Value = 1827
$ws.Range("M141").Value = 3353
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11907907
$ws.Range("I32").Value = 14288251
$ws.Range("J32").Value = 6188.2856
$ws.Range("K32").Value = 14288251
$ws.Range("L32").Value = 6188.2856
$ws.Range("M32").Value = -14287964
$ws.Range("N32").Value = -6762.2856
$ws.Range("H61").Value = 1078868
$ws.Range("I61").Value = 1668855.8
$ws.Range("J61").Value = 6163
$ws.Range("K61").Value = 1668855.8
$ws.Range("L61").Value = 6163
$ws.Range("M61").Value = -1668643.8
$ws.Range("N61").Value = -6587
$ws.Range("H74").Value = 3050967.2
$ws.Range("I74").Value = 3572707
$ws.Range("J74").Value = 7486.5
$ws.Range("K74").Value = 3572707
$ws.Range("L74").Value = 7486.5
$ws.Range("M74").Value = -3571833
$ws.Range("N74").Value = -9234.5
$ws.Range("H77").Value = 3050967.2
$ws.Range("I77").Value = 3572707
$ws.Range("J77").Value = 7486.5
$ws.Range("K77").Value = 17863535
$ws.Range("L77").Value = 37432.5
$ws.Range("M77").Value = -17859167
$ws.Range("N77").Value = -46168.5
$ws.Range("H122").Value = 2616.4666
$ws.Range("I122").Value = 2148.9565
$ws.Range("K122").Value = 6446.869499999999
$ws.Range("M122").Value = -3996.869499999999
$ws.Range("H132").Value = 898031.5
$ws.Range("I132").Value = 1014026.94
$ws.Range("J132").Value = 8733
$ws.Range("K132").Value = 3042080.82
$ws.Range("L132").Value = 26199
$ws.Range("M132").Value = -3039550.82
$ws.Range("N132").Value = -31259
$ws.Range("H136").Value = 1078868
$ws.Range("I136").Value = 1668855.8
$ws.Range("J136").Value = 6163
$ws.Range("K136").Value = 5006567.4
$ws.Range("L136").Value = 18489
$ws.Range("M136").Value = -5004017.4
$ws.Range("N136").Value = -23589
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 23803.607
$ws.Range("I99").Value = 23291.25
$ws.Range("K99").Value = 23291.25
$ws.Range("M99").Value = -21793.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1543.8667
$ws.Range("I22").Value = 669.85
$ws.Range("J22").Value = 3291.9
$ws.Range("K22").Value = 669.85
$ws.Range("L22").Value = 3291.9
$ws.Range("M22").Value = -319.85
$ws.Range("N22").Value = -3991.9
$ws.Range("H31").Value = 96278.69
$ws.Range("I31").Value = 152045.95
$ws.Range("J31").Value = 27389.705
$ws.Range("K31").Value = 152045.95
$ws.Range("L31").Value = 27389.705
$ws.Range("M31").Value = -151750.95
$ws.Range("N31").Value = -27979.705
$ws.Range("H34").Value = 96278.69
$ws.Range("I34").Value = 152045.95
$ws.Range("J34").Value = 27389.705
$ws.Range("K34").Value = 152045.95
$ws.Range("L34").Value = 27389.705
$ws.Range("M34").Value = -151843.95
$ws.Range("N34").Value = -27793.705
$ws.Range("H58").Value = 159999.97
$ws.Range("I58").Value = 229958.53
$ws.Range("J58").Value = 2593.2083
$ws.Range("K58").Value = 229958.53
$ws.Range("L58").Value = 2593.2083
$ws.Range("M58").Value = -229755.53
$ws.Range("N58").Value = -2999.2083
$ws.Range("H80").Value = 47996.5
$ws.Range("J80").Value = 47996.5
$ws.Range("L80").Value = 47996.5
$ws.Range("N80").Value = -50242.5
$ws.Range("H83").Value = 47996.5
$ws.Range("J83").Value = 47996.5
$ws.Range("L83").Value = 143989.5
$ws.Range("N83").Value = -155221.5
$ws.Range("H122").Value = 1824.7188
$ws.Range("I122").Value = 1279.64
$ws.Range("J122").Value = 3771.4285
$ws.Range("K122").Value = 3838.92
$ws.Range("L122").Value = 11314.2855
$ws.Range("M122").Value = -1388.92
$ws.Range("N122").Value = -16214.2855
$ws.Range("H132").Value = 9274135
$ws.Range("I132").Value = 20869.277
$ws.Range("J132").Value = 27780668
$ws.Range("K132").Value = 62607.83099999999
$ws.Range("L132").Value = 83342004
$ws.Range("M132").Value = -60077.83099999999
$ws.Range("N132").Value = -83347064
$ws.Range("H134").Value = 7129.339
$ws.Range("I134").Value = 11033.686
$ws.Range("K134").Value = 33101.058
$ws.Range("M134").Value = -30566.058
$ws.Range("H136").Value = 159999.97
$ws.Range("I136").Value = 229958.53
$ws.Range("J136").Value = 2593.2083
$ws.Range("K136").Value = 689875.59
$ws.Range("L136").Value = 7779.624899999999
$ws.Range("M136").Value = -687325.59
$ws.Range("N136").Value = -12879.6249
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 166796.08
$ws.Range("I14").Value = 166796.08
$ws.Range("K14").Value = 500388.24
$ws.Range("M14").Value = -500215.24
$ws.Range("H17").Value = 58828156
$ws.Range("I17").Value = 56.333332
$ws.Range("J17").Value = 200015600
$ws.Range("K17").Value = 168.999996
$ws.Range("L17").Value = 600046800
$ws.Range("M17").Value = 0.000003999999989900971
$ws.Range("N17").Value = -600047138
$ws.Range("H107").Value = 27778126
$ws.Range("J107").Value = 62500492
$ws.Range("L107").Value = 187501476
$ws.Range("N107").Value = -187505316
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6133.375
$ws.Range("I70").Value = 8839.5
$ws.Range("J70").Value = 3427.25
$ws.Range("K70").Value = 8839.5
$ws.Range("L70").Value = 3427.25
$ws.Range("M70").Value = -8569.5
$ws.Range("N70").Value = -3967.25
$ws.Range("H73").Value = 6133.375
$ws.Range("I73").Value = 8839.5
$ws.Range("J73").Value = 3427.25
$ws.Range("K73").Value = 8839.5
$ws.Range("L73").Value = 3427.25
$ws.Range("M73").Value = -7903.5
$ws.Range("N73").Value = -5299.25
$ws.Range("H132").Value = 1340401.5
$ws.Range("I132").Value = 2009520.5
$ws.Range("J132").Value = 2163.3333
$ws.Range("K132").Value = 6028561.5
$ws.Range("L132").Value = 6489.999899999999
$ws.Range("M132").Value = -6026031.5
$ws.Range("N132").Value = -11549.9999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2961.9412
$ws.Range("I7").Value = 3012.8333
$ws.Range("J7").Value = 2839.8
$ws.Range("K7").Value = 3012.8333
$ws.Range("L7").Value = 2839.8
$ws.Range("M7").Value = -2900.8333
$ws.Range("N7").Value = -3063.8
$ws.Range("H61").Value = 4373.6
$ws.Range("I61").Value = 2167.2856
$ws.Range("J61").Value = 5561.615
$ws.Range("K61").Value = 2167.2856
$ws.Range("L61").Value = 5561.615
$ws.Range("M61").Value = -1965.2856
$ws.Range("N61").Value = -5965.615
$ws.Range("H113").Value = 4373.6
$ws.Range("I113").Value = 2167.2856
$ws.Range("J113").Value = 5561.615
$ws.Range("K113").Value = 2167.2856
$ws.Range("L113").Value = 5561.615
$ws.Range("M113").Value = 2.714399999999841
$ws.Range("N113").Value = -9901.615
$ws.Range("H122").Value = 2970.0312
$ws.Range("I122").Value = 2672.25
$ws.Range("J122").Value = 3863.375
$ws.Range("K122").Value = 8016.75
$ws.Range("L122").Value = 11590.125
$ws.Range("M122").Value = -5566.75
$ws.Range("N122").Value = -16490.125
$ws.Range("H126").Value = 2961.9412
$ws.Range("I126").Value = 3012.8333
$ws.Range("J126").Value = 2839.8
$ws.Range("K126").Value = 9038.499899999999
$ws.Range("L126").Value = 8519.400000000001
$ws.Range("M126").Value = -6568.499899999999
$ws.Range("N126").Value = -13459.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1372.279
$ws.Range("I113").Value = 654.7143
$ws.Range("J113").Value = 2711.7334
$ws.Range("K113").Value = 1964.1429
$ws.Range("L113").Value = 8135.2002
$ws.Range("M113").Value = 205.8571000000002
$ws.Range("N113").Value = -12475.2002
$ws.Range("H132").Value = 5299541.5
$ws.Range("I132").Value = 8051211
$ws.Range("J132").Value = 7869.077
$ws.Range("K132").Value = 24153633
$ws.Range("L132").Value = 23607.231
$ws.Range("M132").Value = -24151103
$ws.Range("N132").Value = -28667.231
